# Add two new "wallet info" history rows (new JSON date entries) below the
# existing data, mirroring the existing A-column date strings.
#
# NOTE: Assigning plain "yyyy-mm-dd" text straight to .Value/.Value2 makes
# Excel auto-recognize it as a date and store it as a numeric date serial
# (with a new date-formatted style) instead of as text. To keep these new
# cells as plain shared-string text - consistent with the existing A2 cell -
# we stage them as literal-string formulas and then convert the formulas to
# their cached values in place via copy / paste-special-values. That leaves
# plain text cells with no formula and no extra cell style, same as the
# existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Formula = "=""2023-08-01"""
$ws.Range("A4").Formula = "=""2023-07-31"""

$ws.Range("A3:A4").Copy()
$ws.Range("A3:A4").PasteSpecial(-4163)  # xlPasteValues
